# Update the "第二学年" (Second Academic Year) bill-details sheet:
# Add two new expense entries for 2019-01-09 (生活费 400) and
# 2019-01-16 (其它/期末聚餐费 200) into the first two empty rows
# of the detail table (rows 24 and 25), then move the active
# selection to F25 as the final editing position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("第二学年")

# Row 24: 2019-01-09, 支出, 400, 生活费, 生活费(2019-01-10 到 放寒假）
$ws.Range("B24").Value = 21
$ws.Range("C24").Value = "支出"
$ws.Range("D24").Value = 400
$ws.Range("E24").Value = 43474
$ws.Range("F24").Value = "生活费"
$ws.Range("G24").Value = "生活费(2019-01-10 到 放寒假）"

# Row 25: 2019-01-16, 支出, 200, 其它, 放假聚餐
$ws.Range("B25").Value = 22
$ws.Range("C25").Value = "支出"
$ws.Range("D25").Value = 200
$ws.Range("E25").Value = 43481
$ws.Range("F25").Value = "其它"
$ws.Range("G25").Value = "放假聚餐"

# Recalculate all formulas (totals / SUMIFS) so the summary boxes
# at the top of the sheet reflect the newly added rows.
$excel.CalculateFull()

# Leave the active cell where the author finished editing.
$ws.Range("F25").Select()
